$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '25.822.27'
$ws.Range('E2').Value = '  -2.60%  '

# Row 3
$ws.Range('D3').Value = '1.640.35'
$ws.Range('E3').Value = '  -1.85%  '

# Row 4
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.020'
$ws.Range('E4').Value = '  +1.22%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '217.09'
$ws.Range('E5').Value = '  -1.10%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.5064'
$ws.Range('E6').Value = '  -2.00%  '

# Row 7
$ws.Range('E7').Value = '  +1.21%  '

# Row 8
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2583'
$ws.Range('E8').Value = '  +0.12%  '

# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06421'
$ws.Range('E9').Value = '  -0.83%  '

# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '19.50'
$ws.Range('E10').Value = '  -2.50%  '

# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07782'
$ws.Range('E11').Value = '  +1.28%  '

# Row 12
$ws.Range('B12').Value = 'Polkadot'
$ws.Range('C12').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '4.262'
$ws.Range('E12').Value = '  -2.11%  '

# Row 13
$ws.Range('B13').Value = 'WrappedEther'
$ws.Range('C13').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D13').Value = '1.636.67'
$ws.Range('E13').Value = '  -2.01%  '

# Row 14
$ws.Range('D14').Value = '1.863.90'

# Row 15
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.5473'
$ws.Range('E15').Value = '  -1.99%  '

# Row 16
$ws.Range('D16').Value = '0.0₅7955'
$ws.Range('E16').Value = '  -1.00%  '

# Row 17
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '63.65'
$ws.Range('E17').Value = '  -1.95%  '

# Row 18
$ws.Range('D18').Value = '25.997.46'
$ws.Range('E18').Value = '  -2.09%  '

# Row 19
$ws.Range('E19').Value = '  +1.01%  '

# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '204.88'
$ws.Range('E20').Value = '  -2.90%  '

# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '4.336'
$ws.Range('E21').Value = '  -2.51%  '

# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '10.01'
$ws.Range('E22').Value = '  -0.95%  '

# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.977'
$ws.Range('E23').Value = '  +1.33%  '

# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.022'
$ws.Range('E24').Value = '  +1.38%  '

# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '1.979'
$ws.Range('E25').Value = '  +15.14%  '

# Row 26
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '142.49'
$ws.Range('E26').Value = '  -0.55%  '

# Row 27
$ws.Range('E27').Value = '  -1.35%  '

# Row 28
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '15.72'
$ws.Range('E28').Value = '  -0.45%  '

# Row 29
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '6.818'
$ws.Range('E29').Value = '  -2.66%  '

# Row 30
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.250'
$ws.Range('E30').Value = '  -1.30%  '

# Row 31
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.05021'
$ws.Range('E31').Value = '  -3.95%  '

# Row 32
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.269'
$ws.Range('E32').Value = '  -2.55%  '

# Row 33
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.210'
$ws.Range('E33').Value = '  -0.23%  '

# Row 34
$ws.Range('E34').Value = '  -3.05%  '

# Row 35
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.365'
$ws.Range('E35').Value = '  -0.60%  '

# Row 36
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.694'
$ws.Range('E36').Value = '  -2.48%  '

# Row 37
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.8942'
$ws.Range('E37').Value = '  -3.34%  '

# Row 38
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.5672'
$ws.Range('E38').Value = '  -1.15%  '

# Row 39
$ws.Range('D39').Value = '1.119.75'
$ws.Range('E39').Value = '  -3.60%  '

# Row 40
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.01568'
$ws.Range('E40').Value = '  -1.64%  '

# Row 41
$ws.Range('B41').Value = 'PaxDollar'
$ws.Range('C41').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.025'
$ws.Range('E41').Value = '  +1.69%  '

# Row 42
$ws.Range('B42').Value = 'mCoin'
$ws.Range('C42').Value = 'https://coinranking.com/coin/fzVgyjBcRc9+mcoin-mcoin'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.597'
$ws.Range('E42').Value = '  +1.21%  '

# Row 43
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '5.661'
$ws.Range('E43').Value = '  +0.25%  '

# Row 44
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.8167'
$ws.Range('E44').Value = '  -2.69%  '

# Row 45
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '100.19'
$ws.Range('E45').Value = '  +0.03%  '

# Row 46
$ws.Range('D46').Value = '1.773.76'
$ws.Range('E46').Value = '  -2.01%  '

# Row 47
$ws.Range('D47').Value = '0.0₈112'
$ws.Range('E47').Value = '  +0.23%  '

# Row 48
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.4560'
$ws.Range('E48').Value = '  +1.43%  '

# Row 49
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.019'
$ws.Range('E49').Value = '  +1.30%  '

# Row 50
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '54.81'
$ws.Range('E50').Value = '  -1.88%  '

# Row 51
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.05044'
$ws.Range('E51').Value = '  -1.81%  '
